$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column C ("Descripcion"), shifting the existing
# "Telefono" column (old C) one place right to become column D.
$ws.Columns("C:C").Insert()
$ws.Range("C1").Value = "Descripcion"

# Fill in the "Descripcion" values for the existing rows (2-9).
$descriptions = @("a", "b", "c", "d", "e", "f", "g", "h")
for ($i = 0; $i -lt $descriptions.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $descriptions[$i]
}

# The old row 9 phone value was stored as text "65"; make it numeric.
$ws.Cells.Item(9, 4).Value = 65

# New client rows (10-17): ID, Nombre, Descripcion, Telefono
$newRows = @(
    @(9,  "lola",    "",        3333),
    @(10, "lolas",   "",        3333),
    @(11, "lolasa",  "dolores", 3333),
    @(12, "lolasa2", "dolores", 3333),
    @(13, "lolasa3", "dolores", 3333),
    @(14, "lolasa4", "dolores", 333333333),
    @(15, "lolasa5", "dolores", 3333333335),
    @(16, "sagra",   "",        "")
)

$r = 10
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
